$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 78
$ws.Range("A78").Value = 7
$ws.Range("B78").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C78").Value = 'Ñuble'
$ws.Range("D78").Value = 45120
$ws.Range("D78").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E78").Value = 16
$ws.Range("F78").Value = 100112013
$ws.Range("G78").Value = 'Alcachofa'
$ws.Range("H78").Value = 'Madrigal'
$ws.Range("I78").Value = 'Primera'
$ws.Range("J78").Value = 40
$ws.Range("K78").Value = 17000
$ws.Range("L78").Value = 17000
$ws.Range("M78").Value = 17000
$ws.Range("N78").Value = '$/caja 40 unidades'
$ws.Range("O78").Value = 'Provincia de Limarí'
$ws.Range("P78").Value = 425
$ws.Range("Q78").Value = 40
$ws.Range("R78").Value = 'Hortaliza'

# Row 79
$ws.Range("A79").Value = 7
$ws.Range("B79").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C79").Value = 'Ñuble'
$ws.Range("D79").Value = 44831
$ws.Range("D79").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E79").Value = 16
$ws.Range("F79").Value = 100112013
$ws.Range("G79").Value = 'Alcachofa'
$ws.Range("H79").Value = 'Española'
$ws.Range("I79").Value = 'Primera'
$ws.Range("J79").Value = 60
$ws.Range("K79").Value = 11000
$ws.Range("L79").Value = 12000
$ws.Range("M79").Value = 11500
$ws.Range("N79").Value = '$/caja 30 unidades'
$ws.Range("O79").Value = 'Provincia de Limarí'
$ws.Range("P79").Value = 383
$ws.Range("Q79").Value = 30
$ws.Range("R79").Value = 'Hortaliza'

# Row 80
$ws.Range("A80").Value = 7
$ws.Range("B80").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C80").Value = 'Ñuble'
$ws.Range("D80").Value = 44831
$ws.Range("D80").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E80").Value = 16
$ws.Range("F80").Value = 100112013
$ws.Range("G80").Value = 'Alcachofa'
$ws.Range("H80").Value = 'Madrigal'
$ws.Range("I80").Value = 'Primera'
$ws.Range("J80").Value = 60
$ws.Range("K80").Value = 10000
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = 10000
$ws.Range("N80").Value = '$/caja 40 unidades'
$ws.Range("O80").Value = 'Provincia de Limarí'
$ws.Range("P80").Value = 250
$ws.Range("Q80").Value = 40
$ws.Range("R80").Value = 'Hortaliza'

# Row 81
$ws.Range("A81").Value = 7
$ws.Range("B81").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C81").Value = 'Ñuble'
$ws.Range("D81").Value = 44473
$ws.Range("D81").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E81").Value = 16
$ws.Range("F81").Value = 100112013
$ws.Range("G81").Value = 'Alcachofa'
$ws.Range("H81").Value = 'Madrigal'
$ws.Range("I81").Value = 'Primera'
$ws.Range("J81").Value = 160
$ws.Range("K81").Value = 11000
$ws.Range("L81").Value = 12000
$ws.Range("M81").Value = 11500
$ws.Range("N81").Value = '$/caja 40 unidades'
$ws.Range("O81").Value = 'Provincia del Elquí'
$ws.Range("P81").Value = 288
$ws.Range("Q81").Value = 40
$ws.Range("R81").Value = 'Hortaliza'

# Row 82
$ws.Range("A82").Value = 7
$ws.Range("B82").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C82").Value = 'Ñuble'
$ws.Range("D82").Value = 44427
$ws.Range("D82").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E82").Value = 16
$ws.Range("F82").Value = 100112013
$ws.Range("G82").Value = 'Alcachofa'
$ws.Range("H82").Value = 'Madrigal'
$ws.Range("I82").Value = 'Primera'
$ws.Range("J82").Value = 120
$ws.Range("K82").Value = 13000
$ws.Range("L82").Value = 14000
$ws.Range("M82").Value = 13500
$ws.Range("N82").Value = '$/caja 40 unidades'
$ws.Range("O82").Value = 'Provincia del Elquí'
$ws.Range("P82").Value = 338
$ws.Range("Q82").Value = 40
$ws.Range("R82").Value = 'Hortaliza'

# Row 83
$ws.Range("A83").Value = 7
$ws.Range("B83").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C83").Value = 'Ñuble'
$ws.Range("D83").Value = 44503
$ws.Range("D83").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E83").Value = 16
$ws.Range("F83").Value = 100112013
$ws.Range("G83").Value = 'Alcachofa'
$ws.Range("H83").Value = 'Madrigal'
$ws.Range("I83").Value = 'Primera'
$ws.Range("J83").Value = 160
$ws.Range("K83").Value = 11000
$ws.Range("L83").Value = 12000
$ws.Range("M83").Value = 11500
$ws.Range("N83").Value = '$/caja 40 unidades'
$ws.Range("O83").Value = 'Provincia del Elquí'
$ws.Range("P83").Value = 288
$ws.Range("Q83").Value = 40
$ws.Range("R83").Value = 'Hortaliza'

# Row 84
$ws.Range("A84").Value = 7
$ws.Range("B84").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C84").Value = 'Ñuble'
$ws.Range("D84").Value = 45097
$ws.Range("D84").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E84").Value = 16
$ws.Range("F84").Value = 100112013
$ws.Range("G84").Value = 'Alcachofa'
$ws.Range("H84").Value = 'Española'
$ws.Range("I84").Value = 'Primera'
$ws.Range("J84").Value = 120
$ws.Range("K84").Value = 15000
$ws.Range("L84").Value = 16000
$ws.Range("M84").Value = 15500
$ws.Range("N84").Value = '$/caja 30 unidades'
$ws.Range("O84").Value = 'Provincia de Limarí'
$ws.Range("P84").Value = 517
$ws.Range("Q84").Value = 30
$ws.Range("R84").Value = 'Hortaliza'

# Row 85
$ws.Range("A85").Value = 7
$ws.Range("B85").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C85").Value = 'Ñuble'
$ws.Range("D85").Value = 44784
$ws.Range("D85").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E85").Value = 16
$ws.Range("F85").Value = 100112013
$ws.Range("G85").Value = 'Alcachofa'
$ws.Range("H85").Value = 'Argentina(o)'
$ws.Range("I85").Value = 'Primera'
$ws.Range("J85").Value = 60
$ws.Range("K85").Value = 14000
$ws.Range("L85").Value = 15000
$ws.Range("M85").Value = 14500
$ws.Range("N85").Value = '$/caja 50 unidades'
$ws.Range("O85").Value = 'Provincia de Limarí'
$ws.Range("P85").Value = 290
$ws.Range("Q85").Value = 50
$ws.Range("R85").Value = 'Hortaliza'

# Row 86
$ws.Range("A86").Value = 7
$ws.Range("B86").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C86").Value = 'Ñuble'
$ws.Range("D86").Value = 44784
$ws.Range("D86").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E86").Value = 16
$ws.Range("F86").Value = 100112013
$ws.Range("G86").Value = 'Alcachofa'
$ws.Range("H86").Value = 'Madrigal'
$ws.Range("I86").Value = 'Primera'
$ws.Range("J86").Value = 60
$ws.Range("K86").Value = 13000
$ws.Range("L86").Value = 14000
$ws.Range("M86").Value = 13500
$ws.Range("N86").Value = '$/caja 40 unidades'
$ws.Range("O86").Value = 'Provincia de Limarí'
$ws.Range("P86").Value = 338
$ws.Range("Q86").Value = 40
$ws.Range("R86").Value = 'Hortaliza'

# Row 87
$ws.Range("A87").Value = 7
$ws.Range("B87").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C87").Value = 'Ñuble'
$ws.Range("D87").Value = 45085
$ws.Range("D87").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E87").Value = 16
$ws.Range("F87").Value = 100112013
$ws.Range("G87").Value = 'Alcachofa'
$ws.Range("H87").Value = 'Española'
$ws.Range("I87").Value = 'Primera'
$ws.Range("J87").Value = 20
$ws.Range("K87").Value = 16000
$ws.Range("L87").Value = 16000
$ws.Range("M87").Value = 16000
$ws.Range("N87").Value = '$/caja 50 unidades'
$ws.Range("O87").Value = 'Provincia de Limarí'
$ws.Range("P87").Value = 320
$ws.Range("Q87").Value = 50
$ws.Range("R87").Value = 'Hortaliza'

# Row 88
$ws.Range("A88").Value = 7
$ws.Range("B88").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C88").Value = 'Ñuble'
$ws.Range("D88").Value = 44838
$ws.Range("D88").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E88").Value = 16
$ws.Range("F88").Value = 100112013
$ws.Range("G88").Value = 'Alcachofa'
$ws.Range("H88").Value = 'Argentina(o)'
$ws.Range("I88").Value = 'Primera'
$ws.Range("J88").Value = 100
$ws.Range("K88").Value = 9000
$ws.Range("L88").Value = 10000
$ws.Range("M88").Value = 9500
$ws.Range("N88").Value = '$/caja 50 unidades'
$ws.Range("O88").Value = 'Provincia de Limarí'
$ws.Range("P88").Value = 190
$ws.Range("Q88").Value = 50
$ws.Range("R88").Value = 'Hortaliza'

# Row 89
$ws.Range("A89").Value = 7
$ws.Range("B89").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C89").Value = 'Ñuble'
$ws.Range("D89").Value = 44838
$ws.Range("D89").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E89").Value = 16
$ws.Range("F89").Value = 100112013
$ws.Range("G89").Value = 'Alcachofa'
$ws.Range("H89").Value = 'Española'
$ws.Range("I89").Value = 'Primera'
$ws.Range("J89").Value = 100
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 11000
$ws.Range("M89").Value = 10500
$ws.Range("N89").Value = '$/caja 30 unidades'
$ws.Range("O89").Value = 'Provincia de Limarí'
$ws.Range("P89").Value = 350
$ws.Range("Q89").Value = 30
$ws.Range("R89").Value = 'Hortaliza'

# Row 90
$ws.Range("A90").Value = 7
$ws.Range("B90").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C90").Value = 'Ñuble'
$ws.Range("D90").Value = 44838
$ws.Range("D90").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E90").Value = 16
$ws.Range("F90").Value = 100112013
$ws.Range("G90").Value = 'Alcachofa'
$ws.Range("H90").Value = 'Madrigal'
$ws.Range("I90").Value = 'Primera'
$ws.Range("J90").Value = 100
$ws.Range("K90").Value = 9000
$ws.Range("L90").Value = 10000
$ws.Range("M90").Value = 9500
$ws.Range("N90").Value = '$/caja 40 unidades'
$ws.Range("O90").Value = 'Provincia de Limarí'
$ws.Range("P90").Value = 238
$ws.Range("Q90").Value = 40
$ws.Range("R90").Value = 'Hortaliza'

# Row 91
$ws.Range("A91").Value = 7
$ws.Range("B91").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C91").Value = 'Ñuble'
$ws.Range("D91").Value = 45079
$ws.Range("D91").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E91").Value = 16
$ws.Range("F91").Value = 100112013
$ws.Range("G91").Value = 'Alcachofa'
$ws.Range("H91").Value = 'Española'
$ws.Range("I91").Value = 'Primera'
$ws.Range("J91").Value = 30
$ws.Range("K91").Value = 16000
$ws.Range("L91").Value = 16000
$ws.Range("M91").Value = 16000
$ws.Range("N91").Value = '$/caja 30 unidades'
$ws.Range("O91").Value = 'Provincia de Limarí'
$ws.Range("P91").Value = 533
$ws.Range("Q91").Value = 30
$ws.Range("R91").Value = 'Hortaliza'

# Row 92
$ws.Range("A92").Value = 7
$ws.Range("B92").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C92").Value = 'Ñuble'
$ws.Range("D92").Value = 44490
$ws.Range("D92").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E92").Value = 16
$ws.Range("F92").Value = 100112013
$ws.Range("G92").Value = 'Alcachofa'
$ws.Range("H92").Value = 'Madrigal'
$ws.Range("I92").Value = 'Primera'
$ws.Range("J92").Value = 100
$ws.Range("K92").Value = 11000
$ws.Range("L92").Value = 12000
$ws.Range("M92").Value = 11500
$ws.Range("N92").Value = '$/caja 40 unidades'
$ws.Range("O92").Value = 'Provincia del Elquí'
$ws.Range("P92").Value = 288
$ws.Range("Q92").Value = 40
$ws.Range("R92").Value = 'Hortaliza'

# Row 93
$ws.Range("A93").Value = 7
$ws.Range("B93").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C93").Value = 'Ñuble'
$ws.Range("D93").Value = 44495
$ws.Range("D93").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E93").Value = 16
$ws.Range("F93").Value = 100112013
$ws.Range("G93").Value = 'Alcachofa'
$ws.Range("H93").Value = 'Madrigal'
$ws.Range("I93").Value = 'Primera'
$ws.Range("J93").Value = 120
$ws.Range("K93").Value = 11000
$ws.Range("L93").Value = 12000
$ws.Range("M93").Value = 11500
$ws.Range("N93").Value = '$/caja 40 unidades'
$ws.Range("O93").Value = 'Provincia del Elquí'
$ws.Range("P93").Value = 288
$ws.Range("Q93").Value = 40
$ws.Range("R93").Value = 'Hortaliza'

# Row 94
$ws.Range("A94").Value = 7
$ws.Range("B94").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C94").Value = 'Ñuble'
$ws.Range("D94").Value = 44417
$ws.Range("D94").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E94").Value = 16
$ws.Range("F94").Value = 100112013
$ws.Range("G94").Value = 'Alcachofa'
$ws.Range("H94").Value = 'Madrigal'
$ws.Range("I94").Value = 'Primera'
$ws.Range("J94").Value = 120
$ws.Range("K94").Value = 15000
$ws.Range("L94").Value = 16000
$ws.Range("M94").Value = 15500
$ws.Range("N94").Value = '$/caja 40 unidades'
$ws.Range("O94").Value = 'Provincia del Elquí'
$ws.Range("P94").Value = 388
$ws.Range("Q94").Value = 40
$ws.Range("R94").Value = 'Hortaliza'

# Row 95
$ws.Range("A95").Value = 7
$ws.Range("B95").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C95").Value = 'Ñuble'
$ws.Range("D95").Value = 44468
$ws.Range("D95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E95").Value = 16
$ws.Range("F95").Value = 100112013
$ws.Range("G95").Value = 'Alcachofa'
$ws.Range("H95").Value = 'Madrigal'
$ws.Range("I95").Value = 'Primera'
$ws.Range("J95").Value = 60
$ws.Range("K95").Value = 12000
$ws.Range("L95").Value = 13000
$ws.Range("M95").Value = 12500
$ws.Range("N95").Value = '$/caja 40 unidades'
$ws.Range("O95").Value = 'Provincia del Elquí'
$ws.Range("P95").Value = 312
$ws.Range("Q95").Value = 40
$ws.Range("R95").Value = 'Hortaliza'

# Row 96
$ws.Range("A96").Value = 7
$ws.Range("B96").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C96").Value = 'Ñuble'
$ws.Range("D96").Value = 44420
$ws.Range("D96").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E96").Value = 16
$ws.Range("F96").Value = 100112013
$ws.Range("G96").Value = 'Alcachofa'
$ws.Range("H96").Value = 'Madrigal'
$ws.Range("I96").Value = 'Primera'
$ws.Range("J96").Value = 120
$ws.Range("K96").Value = 13000
$ws.Range("L96").Value = 14000
$ws.Range("M96").Value = 13500
$ws.Range("N96").Value = '$/caja 40 unidades'
$ws.Range("O96").Value = 'Provincia del Elquí'
$ws.Range("P96").Value = 338
$ws.Range("Q96").Value = 40
$ws.Range("R96").Value = 'Hortaliza'

# Row 97
$ws.Range("A97").Value = 7
$ws.Range("B97").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C97").Value = 'Ñuble'
$ws.Range("D97").Value = 44811
$ws.Range("D97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E97").Value = 16
$ws.Range("F97").Value = 100112013
$ws.Range("G97").Value = 'Alcachofa'
$ws.Range("H97").Value = 'Española'
$ws.Range("I97").Value = 'Primera'
$ws.Range("J97").Value = 60
$ws.Range("K97").Value = 12000
$ws.Range("L97").Value = 13000
$ws.Range("M97").Value = 12500
$ws.Range("N97").Value = '$/caja 30 unidades'
$ws.Range("O97").Value = 'Provincia de Limarí'
$ws.Range("P97").Value = 417
$ws.Range("Q97").Value = 30
$ws.Range("R97").Value = 'Hortaliza'

# Row 98
$ws.Range("A98").Value = 7
$ws.Range("B98").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C98").Value = 'Ñuble'
$ws.Range("D98").Value = 45112
$ws.Range("D98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E98").Value = 16
$ws.Range("F98").Value = 100112013
$ws.Range("G98").Value = 'Alcachofa'
$ws.Range("H98").Value = 'Argentina(o)'
$ws.Range("I98").Value = 'Primera'
$ws.Range("J98").Value = 50
$ws.Range("K98").Value = 16000
$ws.Range("L98").Value = 16000
$ws.Range("M98").Value = 16000
$ws.Range("N98").Value = '$/caja 50 unidades'
$ws.Range("O98").Value = 'Provincia de Limarí'
$ws.Range("P98").Value = 320
$ws.Range("Q98").Value = 50
$ws.Range("R98").Value = 'Hortaliza'

# Row 99
$ws.Range("A99").Value = 7
$ws.Range("B99").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C99").Value = 'Ñuble'
$ws.Range("D99").Value = 45089
$ws.Range("D99").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E99").Value = 16
$ws.Range("F99").Value = 100112013
$ws.Range("G99").Value = 'Alcachofa'
$ws.Range("H99").Value = 'Argentina(o)'
$ws.Range("I99").Value = 'Primera'
$ws.Range("J99").Value = 60
$ws.Range("K99").Value = 16000
$ws.Range("L99").Value = 17000
$ws.Range("M99").Value = 16500
$ws.Range("N99").Value = '$/caja 50 unidades'
$ws.Range("O99").Value = 'Provincia de Limarí'
$ws.Range("P99").Value = 330
$ws.Range("Q99").Value = 50
$ws.Range("R99").Value = 'Hortaliza'

# Row 100
$ws.Range("A100").Value = 7
$ws.Range("B100").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C100").Value = 'Ñuble'
$ws.Range("D100").Value = 45121
$ws.Range("D100").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E100").Value = 16
$ws.Range("F100").Value = 100112013
$ws.Range("G100").Value = 'Alcachofa'
$ws.Range("H100").Value = 'Madrigal'
$ws.Range("I100").Value = 'Primera'
$ws.Range("J100").Value = 40
$ws.Range("K100").Value = 17000
$ws.Range("L100").Value = 17000
$ws.Range("M100").Value = 17000
$ws.Range("N100").Value = '$/caja 40 unidades'
$ws.Range("O100").Value = 'Provincia de Limarí'
$ws.Range("P100").Value = 425
$ws.Range("Q100").Value = 40
$ws.Range("R100").Value = 'Hortaliza'

# Row 101
$ws.Range("A101").Value = 7
$ws.Range("B101").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C101").Value = 'Ñuble'
$ws.Range("D101").Value = 44777
$ws.Range("D101").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E101").Value = 16
$ws.Range("F101").Value = 100112013
$ws.Range("G101").Value = 'Alcachofa'
$ws.Range("H101").Value = 'Madrigal'
$ws.Range("I101").Value = 'Primera'
$ws.Range("J101").Value = 60
$ws.Range("K101").Value = 14000
$ws.Range("L101").Value = 15000
$ws.Range("M101").Value = 14500
$ws.Range("N101").Value = '$/caja 40 unidades'
$ws.Range("O101").Value = 'Provincia del Elquí'
$ws.Range("P101").Value = 362
$ws.Range("Q101").Value = 40
$ws.Range("R101").Value = 'Hortaliza'
